$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos price/volume snapshot (and the Kaspa/InjectiveProtocol
# and ApeXProtocol/ThetaToken row swaps) to the new scrape values.
# Price-column values that look like plain decimals (e.g. "9.10") are
# written with a leading apostrophe so Excel keeps them as literal text
# instead of auto-converting to a Number and stripping trailing zeros
# (e.g. "9.10" -> 9.1) -- matching the source data's text cell type.
$ws.Range("D2").Value = '66.280.83'
$ws.Range("E2").Value = '  +3.60%  '
$ws.Range("D3").Value = '3.492.76'
$ws.Range("E3").Value = '  +4.56%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''555.92'
$ws.Range("E5").Value = '  +6.59%  '
$ws.Range("D6").Value = '''185.11'
$ws.Range("E6").Value = '  +7.15%  '
$ws.Range("D7").Value = '''0.645'
$ws.Range("E7").Value = '  +9.21%  '
$ws.Range("D8").Value = '3.486.87'
$ws.Range("E8").Value = '  +4.47%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").Value = '''0.633'
$ws.Range("E10").Value = '  +4.70%  '
$ws.Range("E11").Value = '  +15.20%  '
$ws.Range("D12").Value = '''54.55'
$ws.Range("E12").Value = '  +2.99%  '
$ws.Range("D13").Value = '''0.0000271'
$ws.Range("E13").Value = '  +5.91%  '
$ws.Range("D14").Value = '''9.29'
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("D15").Value = '4.063.78'
$ws.Range("E15").Value = '  +4.40%  '
$ws.Range("D16").Value = '3.499.23'
$ws.Range("E16").Value = '  +4.42%  '
$ws.Range("D17").Value = '''18.58'
$ws.Range("E17").Value = '  +6.55%  '
$ws.Range("E18").Value = '  +3.60%  '
$ws.Range("D19").Value = '66.365.22'
$ws.Range("E19").Value = '  +3.95%  '
$ws.Range("D20").Value = '''12.02'
$ws.Range("E20").Value = '  +7.57%  '
$ws.Range("D21").Value = '''0.993'
$ws.Range("E21").Value = '  +3.82%  '
$ws.Range("D22").Value = '''422.06'
$ws.Range("E22").Value = '  +13.14%  '
$ws.Range("D23").Value = '''4.06'
$ws.Range("E23").Value = '  +10.85%  '
$ws.Range("D24").Value = '''86.33'
$ws.Range("E24").Value = '  +6.10%  '
$ws.Range("D25").Value = '''4.14'
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("D26").Value = '''10.93'
$ws.Range("E26").Value = '  -4.16%  '
$ws.Range("D27").Value = '''2.91'
$ws.Range("E27").Value = '  +8.01%  '
$ws.Range("D28").Value = '''12.32'
$ws.Range("E28").Value = '  +9.73%  '
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("D30").Value = '''9.10'
$ws.Range("E30").Value = '  +11.52%  '
$ws.Range("D31").Value = '''30.18'
$ws.Range("E31").Value = '  +5.09%  '
$ws.Range("D32").Value = '''631.18'
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").Value = '''6.60'
$ws.Range("E33").Value = '  +3.25%  '
$ws.Range("D34").Value = '''11.73'
$ws.Range("E34").Value = '  +5.21%  '
$ws.Range("E35").Value = '  +5.45%  '
$ws.Range("D36").Value = '''59.97'
$ws.Range("E36").Value = '  +3.70%  '
$ws.Range("D37").Value = '0.0₃0811'
$ws.Range("E37").Value = '  +11.44%  '
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").Value = '''37.82'
$ws.Range("E38").Value = '  +5.46%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.146'
$ws.Range("E39").Value = '  +18.93%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").Value = '''0.385'
$ws.Range("E41").Value = '  +2.22%  '
$ws.Range("D42").Value = '''3.41'
$ws.Range("E42").Value = '  +12.94%  '
$ws.Range("D43").Value = '3.110.40'
$ws.Range("E43").Value = '  +5.95%  '
$ws.Range("D44").Value = '''0.999'
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '''2.60'
$ws.Range("E45").Value = '  -1.46%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''3.38'
$ws.Range("E46").Value = '  +13.11%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").Value = '''2.84'
$ws.Range("E47").Value = '  +9.86%  '
$ws.Range("D48").Value = '''0.0414'
$ws.Range("E48").Value = '  +5.18%  '
$ws.Range("D49").Value = '''2.72'
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("D50").Value = '''0.135'
$ws.Range("E50").Value = '  +8.35%  '
$ws.Range("D51").Value = '''138.99'
$ws.Range("E51").Value = '  +2.38%  '
